$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) values are written as text, matching the
# original inline-string cell type, so Excel does not coerce numeric-
# looking strings (e.g. "15.20", "1.00") into numbers and strip
# formatting / trailing zeros.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.435.45"
$ws.Range("D3").Value = "2.313.50"
$ws.Range("D5").Value = "310.74"
$ws.Range("D6").Value = "104.73"
$ws.Range("D9").Value = "0.529"
$ws.Range("D10").Value = "36.76"
$ws.Range("D11").Value = "52.88"
$ws.Range("D13").Value = "0.112"
$ws.Range("D14").Value = "7.01"
$ws.Range("D15").Value = "2.672.23"
$ws.Range("D16").Value = "15.20"
$ws.Range("D17").Value = "2.311.96"
$ws.Range("D19").Value = "43.338.76"
$ws.Range("D20").Value = "12.17"
$ws.Range("D21").Value = "0.0₃0925"
$ws.Range("D23").Value = "68.24"
$ws.Range("D24").Value = "242.34"
$ws.Range("D26").Value = "2.61"
$ws.Range("D27").Value = "1.00"
$ws.Range("D30").Value = "37.20"
$ws.Range("D32").Value = "166.89"
$ws.Range("D33").Value = "5.29"
$ws.Range("D35").Value = "18.30"
$ws.Range("D37").Value = "0.0745"
$ws.Range("D38").Value = "3.06"
$ws.Range("D39").Value = "4.52"
$ws.Range("D40").Value = "1.88"
$ws.Range("D43").Value = "2.73"
$ws.Range("D45").Value = "1.992.91"
$ws.Range("D46").Value = "19.07"
$ws.Range("D47").Value = "3.06"
$ws.Range("D48").Value = "10.01"
$ws.Range("D49").Value = "57.01"
$ws.Range("D50").Value = "2.96"

# Restore the default "Normal" style on column D so the cells retain
# no explicit style/number-format override, just like the original file.
$dRange.Style = "Normal"

# Coin name / link swap (rows 39-40) plus other text/value updates.
$ws.Range("B39").Value = "RenderToken"
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"

# Volume(1h) percentage updates.
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +7.18%  "
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +8.12%  "
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("E16").Value = "  +3.98%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +5.27%  "
$ws.Range("E29").Value = "  +12.06%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("E36").Value = "  +6.73%  "
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  +8.25%  "
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("E43").Value = "  +20.54%  "
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("E47").Value = "  +3.84%  "
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("E49").Value = "  +5.14%  "
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("E51").Value = "  +8.70%  "
